# ResonatorArray.xlsx edit: add TVPA / AGA Marks coordinate tables in rows 4-10
# and widen a few columns.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column width tweaks -------------------------------------------------
# Internal stored width (OOXML "width") = round(CharWidth*6)/6 + 5/6, so to
# land exactly on an integer stored width W we need CharWidth = (W*6-5)/6.
$ws.Columns.Item(2).ColumnWidth = 10.166666666666666   # B: 6  -> 11
$ws.Columns.Item(4).ColumnWidth = 7.166666666666667    # D: 7  -> 8
$ws.Columns.Item(7).ColumnWidth = 6.166666666666667    # G: 6  -> 7

# --- Populate rows 4-10 with the same cell style used by the existing
#     table (row 11 is a good, fully-populated exemplar of style index 1)
$ws.Range("A11:Y11").Copy() | Out-Null
$ws.Range("A4:Y10").PasteSpecial(-4122) | Out-Null

# Row 4: section headers
$ws.Range("C4").Value = "AGA Marks"
$ws.Range("H4").Value = "TVPA"

# Row 5: column headers for the two coordinate blocks
$ws.Range("B5").Value = "direction"
$ws.Range("C5").Value = "layer"
$ws.Range("D5").Value = "x"
$ws.Range("E5").Value = "y"
$ws.Range("G5").Value = "layer"
$ws.Range("H5").Value = "x"
$ws.Range("I5").Value = "y"

# Row 6: first data row
$ws.Range("B6").Value = "y"
$ws.Range("C6").Value = "GP"
$ws.Range("D6").Value = -2.875
$ws.Range("E6").Value = 0
$ws.Range("G6").Value = "GP"
$ws.Range("H6").Value = -2.525
$ws.Range("I6").Value = 0

# Row 7: second data row
$ws.Range("B7").Value = "x"
$ws.Range("C7").Value = "GP"
$ws.Range("D7").Value = -2.775
$ws.Range("E7").Value = 0
